$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of item data loaded into the sheet (rows 7-20).
# Column A repeats the existing "Some name." shared string; column B
# holds currency-formatted quantities/prices matching the style already
# used by the existing rows (B2:B6).
$values = @(12, 12, 12, 12, 12, 12, 12, 13, 14, 15, 16, 17, 18, 19)
$currencyFormat = $ws.Cells.Item(2, 2).NumberFormat

$row = 7
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = "Some name."
    $ws.Cells.Item($row, 2).Value = $v
    $ws.Cells.Item($row, 2).NumberFormat = $currencyFormat
    $row = $row + 1
}

# Scroll the view down and select E17, matching where the user ended up
# after loading the data.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E17").Select()
